# Generate Report for Handback
# b578155d-b0bc-4d93-9646-485830bba23d.md and d5e47455-8918-4c58-975d-52ce5e2b11db.md
# have now been handed back for both zh-cn and de-de, so the report needs to
# reflect their new status, target/handback files and handback timestamps.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet - zh-cn (E) / de-de (F) status columns for the two files
# that just got handed back.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = $handedBack
$overview.Range("F4").Value = $handedBack
$overview.Range("E5").Value = $handedBack
$overview.Range("F5").Value = $handedBack

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C4").Value = $handedBack
$zhcn.Range("I4").Value = "b578155d-b0bc-4d93-9646-485830bba23d.76ccf2207294c4db3d60bc637a4d9b7518eed91f.zh-cn.xlf"
$zhcn.Range("J4").Value = "b578155d-b0bc-4d93-9646-485830bba23d.76ccf2207294c4db3d60bc637a4d9b7518eed91f.zh-cn.xlf"
$zhcn.Range("K4").Value = "2016-09-02 02:31:34"

$zhcn.Range("C5").Value = $handedBack
$zhcn.Range("I5").Value = "d5e47455-8918-4c58-975d-52ce5e2b11db.f0395604b4d2a0eb39a91ceb93b2dcb4de4526f2.zh-cn.xlf"
$zhcn.Range("J5").Value = "d5e47455-8918-4c58-975d-52ce5e2b11db.f0395604b4d2a0eb39a91ceb93b2dcb4de4526f2.zh-cn.xlf"
$zhcn.Range("K5").Value = "2016-09-02 02:31:34"

$zhcn.Hyperlinks.Add($zhcn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fc3feefeef6fadb931286a9761b34d8838602d3f/e2e/b578155d-b0bc-4d93-9646-485830bba23d.md", "", "", "b578155d-b0bc-4d93-9646-485830bba23d.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fc3feefeef6fadb931286a9761b34d8838602d3f/e2e/d5e47455-8918-4c58-975d-52ce5e2b11db.md", "", "", "d5e47455-8918-4c58-975d-52ce5e2b11db.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C4").Value = $handedBack
$dede.Range("I4").Value = "b578155d-b0bc-4d93-9646-485830bba23d.76ccf2207294c4db3d60bc637a4d9b7518eed91f.de-de.xlf"
$dede.Range("J4").Value = "b578155d-b0bc-4d93-9646-485830bba23d.76ccf2207294c4db3d60bc637a4d9b7518eed91f.de-de.xlf"
$dede.Range("K4").Value = "2016-09-02 02:31:41"

$dede.Range("C5").Value = $handedBack
$dede.Range("I5").Value = "d5e47455-8918-4c58-975d-52ce5e2b11db.f0395604b4d2a0eb39a91ceb93b2dcb4de4526f2.de-de.xlf"
$dede.Range("J5").Value = "d5e47455-8918-4c58-975d-52ce5e2b11db.f0395604b4d2a0eb39a91ceb93b2dcb4de4526f2.de-de.xlf"
$dede.Range("K5").Value = "2016-09-02 02:31:41"

$dede.Hyperlinks.Add($dede.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/26c8cc1bdc10b95b641c6e458ce46e756144e14b/e2e/b578155d-b0bc-4d93-9646-485830bba23d.md", "", "", "b578155d-b0bc-4d93-9646-485830bba23d.md")
$dede.Hyperlinks.Add($dede.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/26c8cc1bdc10b95b641c6e458ce46e756144e14b/e2e/d5e47455-8918-4c58-975d-52ce5e2b11db.md", "", "", "d5e47455-8918-4c58-975d-52ce5e2b11db.md")
